$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right after
#    the title heading (paragraph #2 in the document).
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert a new paragraph, containing the bold title text, right before the
#    very last paragraph (the one that holds the image-generation prompt).
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore() | Out-Null

# A brand new (still empty) paragraph now exists right before the last one;
# grab its range and fill it with the desired OOXML so that it ends up with
# exactly the same run/formatting shape used elsewhere in the document
# (an empty leading run followed by a bold run), with no stray formatting
# inherited from neighboring paragraphs.
$newParaIndex = $d.Paragraphs.Count - 1
$newPara = $d.Paragraphs($newParaIndex)
$newRange = $newPara.Range

$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruit Super Nova for Free: Online Slot Game Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($newParaXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Replace the text of the final paragraph (the italic image-prompt
#    paragraph) with the new meta-description copy, keeping its formatting.
# ---------------------------------------------------------------------------
$oldPrompt = "Please create a feature image that complements the retro theme of the " + [char]34 + "Fruit Super Nova" + [char]34 + " game. The image should showcase a happy Maya warrior wearing glasses. The image should be in cartoon style, with a colorful and eye-catching design. The warrior can be seen holding a basket of brightly colored fruits, with a flaming star (the scatter symbol) just above the basket. The background can feature a starry night sky, with the silhouette of a tropical forest in the distance. The image should convey a fun and exciting atmosphere, inviting players to try out the game and discover its simplicity and entertainment value."
$newDescription = "Uncover the jackpots associated with fruit in Fruit Super Nova. Play this retro-style slot game for free and enjoy fast and straightforward gameplay."

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2) | Out-Null
